$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6951934695243835
$ws.Range("B1").Value = 2.043120861053467
$ws.Range("C1").Value = 2.454780340194702
$ws.Range("D1").Value = 0.80034339427948
$ws.Range("E1").Value = 0.9028966426849365
